# FedExShipments.xlsx - "changes of 18th june 2022"
# Refresh the ShipmentTracking (P), ActualRate (Q) and Result (R) columns
# with the latest FedEx run results.
#
# P/Q/R cells in this sheet are stored as plain text (no number formatting,
# no explicit cell style). Setting .Value directly on a range makes Excel's
# COM layer auto-detect numeric-looking strings (tracking numbers, currency
# amounts) and coerce them into numbers (picking up a number style in the
# process). To keep them as genuine text cells - matching the original
# workbook layout exactly - we briefly mark the destination range as text
# ("@") before assigning the value, then drop back to the default "Normal"
# style so no stray style id is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# row -> (ShipmentTracking, ActualRate, Result)
Set-TextCell $ws.Range("P2")  "320018511882"
Set-TextCell $ws.Range("Q2")  "`$61.34"
$ws.Range("R2").Value = "FAIL"

Set-TextCell $ws.Range("P3")  "320018511996"
Set-TextCell $ws.Range("Q3")  "`$69.80"
$ws.Range("R3").Value = "FAIL"

Set-TextCell $ws.Range("P4")  "320018512021"
Set-TextCell $ws.Range("Q4")  "`$74.03"
$ws.Range("R4").Value = "FAIL"

Set-TextCell $ws.Range("P5")  "320018512098"
Set-TextCell $ws.Range("Q5")  "`$85.66"
$ws.Range("R5").Value = "FAIL"

Set-TextCell $ws.Range("P6")  "320018512135"
Set-TextCell $ws.Range("Q6")  "`$98.35"
$ws.Range("R6").Value = "FAIL"

Set-TextCell $ws.Range("P7")  "320018512205"
Set-TextCell $ws.Range("Q7")  "`$273.38"
$ws.Range("R7").Value = "FAIL"

Set-TextCell $ws.Range("P8")  "320018512238"
Set-TextCell $ws.Range("Q8")  "`$61.34"
$ws.Range("R8").Value = "FAIL"

Set-TextCell $ws.Range("P9")  "320018512250"
Set-TextCell $ws.Range("Q9")  "`$65.57"
$ws.Range("R9").Value = "FAIL"

Set-TextCell $ws.Range("P10") "320018512282"
Set-TextCell $ws.Range("Q10") "`$69.80"
$ws.Range("R10").Value = "FAIL"

Set-TextCell $ws.Range("P11") "320018512319"
Set-TextCell $ws.Range("Q11") "`$82.49"
$ws.Range("R11").Value = "FAIL"

Set-TextCell $ws.Range("P12") "320018512400"
Set-TextCell $ws.Range("Q12") "`$95.18"
$ws.Range("R12").Value = "FAIL"

# Rows 13-19 and 21-26 only got a new tracking number; ActualRate/Result
# were left as-is.
Set-TextCell $ws.Range("P13") "320018503130"
Set-TextCell $ws.Range("P14") "320018503163"
Set-TextCell $ws.Range("P15") "320018503185"
Set-TextCell $ws.Range("P16") "320018503211"
Set-TextCell $ws.Range("P17") "320018503244"
Set-TextCell $ws.Range("P18") "320018503288"
Set-TextCell $ws.Range("P19") "320018503303"

Set-TextCell $ws.Range("P20") "320018503336"
Set-TextCell $ws.Range("Q20") "`$104.69"
$ws.Range("R20").Value = "FAIL"

Set-TextCell $ws.Range("P21") "320018503358"
Set-TextCell $ws.Range("P22") "320018503380"
Set-TextCell $ws.Range("P23") "320018503391"
Set-TextCell $ws.Range("P24") "320018503406"
Set-TextCell $ws.Range("P25") "320018503417"
Set-TextCell $ws.Range("P26") "320018503428"
